# Phase 4 Use Cases: fix pre-conditions for "Basic Radio Experimentation".
#
# The pre-conditions cell had extraneous material in the middle of its
# sentence ("They have working Phase 4 radios. They have enough technical
# expertise to experiment with using a new codec on their radios."). Remove
# it, which is exactly what happened where the (hidden) _GoBack bookmark
# used to sit at the end of the document -- Word relocates that bookmark to
# track the point of the last edit, so move it here too.

$d = $word.ActiveDocument

$marker = "Alice and Bob are licensed operators. They are registered on the Phase 4 system. They have working Phase 4 radios. They have enough technical expertise to experiment with using a new codec on their radios. They are able to have successful 2-way voice communications. "
$part1  = "Alice and Bob are licensed operators. They are registered on the Phase 4 system."
$part2  = " They are able to have successful 2-way voice communications. "

# Locate the (unique) original sentence.
$found = $d.Content
$ok = $found.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Pre-conditions sentence found:" $ok

if ($ok) {
    $baseStart = $found.Start
    $baseEnd   = $found.End

    # Replace the whole sentence with the trimmed-down version (drop the
    # extraneous middle sentence).
    $whole = $d.Range($baseStart, $baseEnd)
    $whole.Text = ($part1 + $part2)

    $bmPos = $baseStart + $part1.Length

    # Drop the stale _GoBack bookmark (currently sitting in its own
    # paragraph near the end of the document) before re-adding it at the
    # new location.
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }

    $bm = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bm)

    # Splitting the run at the bookmark leaves a spurious
    # xml:space="preserve" on the first (no-longer-padded) fragment;
    # round-trip its text through a placeholder so it gets re-serialized
    # without the now-unneeded attribute.
    $frag1 = $d.Range($baseStart, $bmPos)
    $frag1.Text = "TEMP_PLACEHOLDER"
    $frag1b = $d.Range($baseStart, $baseStart + "TEMP_PLACEHOLDER".Length)
    $frag1b.Text = $part1
}
